$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 new rows before the existing "McNemar vs LogReg" block (old row 5),
# pushing it (and everything after) down from rows 5-9 to rows 17-21.
$ws.Rows("5:16").Insert()

# Copy formats (bold header / plain-int / 0.000-numeric) from the existing
# "Accuracy (Over 4 Folds)" block (rows 1-3) onto each of the three new blocks,
# then fill in the text/numbers.

# --- New block: "Precision (Over 4 Folds)" (rows 5-8) ---
$ws.Range("A1").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
$ws.Range("B3:E3").Copy()
$ws.Range("B8:E8").PasteSpecial(-4122)

$ws.Range("A5").Value = "Precision (Over 4 Folds)"
$ws.Range("A6").Value = "n"
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 18
$ws.Range("A7").Value = "Random Forest"
$ws.Range("B7").Value = 0.95199999999999996
$ws.Range("C7").Value = 0.86609999999999998
$ws.Range("D7").Value = 0.59519999999999995
$ws.Range("E7").Value = 0.50990000000000002

# --- New block: "Recall (Over 4 Folds)" (rows 9-12) ---
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("B3:E3").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)

$ws.Range("A9").Value = "Recall (Over 4 Folds)"
$ws.Range("A10").Value = "n"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 15
$ws.Range("E10").Value = 18
$ws.Range("A11").Value = "Random Forest"
$ws.Range("B11").Value = 0.95189999999999997
$ws.Range("C11").Value = 0.86580000000000001
$ws.Range("D11").Value = 0.59619999999999995
$ws.Range("E11").Value = 0.50970000000000004

# --- New block: "F1 (Over 4 Folds)" (rows 13-16) ---
$ws.Range("A1").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A3:E3").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("B3:E3").Copy()
$ws.Range("B16:E16").PasteSpecial(-4122)

$ws.Range("A13").Value = "F1 (Over 4 Folds)"
$ws.Range("A14").Value = "n"
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 15
$ws.Range("E14").Value = 18
$ws.Range("A15").Value = "Random Forest"
$ws.Range("B15").Value = 0.95189999999999997
$ws.Range("C15").Value = 0.86580000000000001
$ws.Range("D15").Value = 0.5948
$ws.Range("E15").Value = 0.50939999999999996

# Move the selection to match the committed workbook state.
$ws.Range("H10").Select()
